$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename products
$ws.Range("A1").Value = "Chocolate"
$ws.Range("A2").Value = "Crisp Chips"

# Update quantities / prices
$ws.Range("B1").Value = 13243
$ws.Range("C1").Value = 30.5
$ws.Range("B2").Value = 2346
$ws.Range("D2").Value = 18.5

# Turn F2 from a text date into a real date value, matching F1's date format
$ws.Range("F1").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = (Get-Date -Year 2023 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)

# Update selection
$ws.Range("D2").Select()
